# Replace original Worksheet content with new Worksheet content:
# every previously-numeric-typed cell in the data sheet becomes a text
# cell (literal strings instead of mis-typed numeric values), while the
# already-numeric price/weight cells are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -eq "") {
        $cell.Value = "'"
    } else {
        $cell.Value = "'" + $text
    }
    # Drop the quote-prefix style Excel auto-applies so the cell keeps the
    # default (unstyled) format, matching the source worksheet.
    $cell.Style = "Normal"
}

# Row 1 - A1/B1 stay blank (no real content); C1/D1/E1 hold literal text
# (previously these were bogus numeric <v> payloads like str('Unknown')).
Set-TextCell 1 1 ""
Set-TextCell 1 2 ""
Set-TextCell 1 3 "str('Unknown')"
Set-TextCell 1 4 "float(4.5)"
Set-TextCell 1 5 "int(500)"

# Row 2 - header labels as text
Set-TextCell 2 1 "id"
Set-TextCell 2 2 "productName"
Set-TextCell 2 3 "productType"
Set-TextCell 2 4 "price"
Set-TextCell 2 5 "weight"

# Row 3 - product 001 pork meat (text) ; price/weight remain numeric
Set-TextCell 3 1 "001"
Set-TextCell 3 2 "pork"
Set-TextCell 3 3 "meat"
$ws.Cells.Item(3, 4).Value = 2.5
$ws.Cells.Item(3, 5).Value = 1000

# Row 4 - product 002 beef meat (text) ; price/weight remain numeric
Set-TextCell 4 1 "002"
Set-TextCell 4 2 "beef"
Set-TextCell 4 3 "meat"
$ws.Cells.Item(4, 4).Value = 4.5
$ws.Cells.Item(4, 5).Value = 1000
